$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.678679804978435
$ws.Range("K2").Value = 0.67987520077397
$ws.Range("L2").Value = 0.74084626345024
$ws.Range("N2").Value = 0.640735542026883
$ws.Range("B3").Value = 0.553718303775116
$ws.Range("D3").Value = 0.62468719027894
$ws.Range("E3").Value = 0.569498424923018
$ws.Range("F3").Value = 0.759057679811368
$ws.Range("G3").Value = 0.522861429260386
$ws.Range("H3").Value = 0.68963033430472
$ws.Range("I3").Value = 0.593566149815293
$ws.Range("J3").Value = 0.610091501584142
$ws.Range("K3").Value = 0.543669469571705
$ws.Range("L3").Value = 0.489608113420232
$ws.Range("M3").Value = 0.845257807486015
$ws.Range("N3").Value = 0.476878037078152
$ws.Range("B4").Value = 0.647240480131028
$ws.Range("K4").Value = 0.665933688883394
$ws.Range("L4").Value = 0.572690723237937
$ws.Range("N4").Value = 0.603471262597791
$ws.Range("B5").Value = 0.672054760884641
$ws.Range("C5").Value = 0.73150166637206
$ws.Range("K5").Value = 0.67292535125419
$ws.Range("L5").Value = 0.639033952437169
$ws.Range("N5").Value = 0.589468339412676
$ws.Range("B6").Value = 0.68518253119246
$ws.Range("K6").Value = 0.708149579691965
$ws.Range("L6").Value = 0.604866837966034
$ws.Range("N6").Value = 0.636828257869627
$ws.Range("B7").Value = 0.61895417279981
$ws.Range("K7").Value = 0.633861854984241
$ws.Range("L7").Value = 0.5411746546684
$ws.Range("N7").Value = 0.560767240057917
